# admin access with dataSourceResource
# Mark several "feature" rows as done, and flag a known issue on the
# "view order" row with a highlighted note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feature")

# Status column (E) -> "done" for checkout/contact-info, credit-info,
# order/placeorder and invoice rows.
$ws.Range("E8").Value  = "done"
$ws.Range("E9").Value  = "done"
$ws.Range("E11").Value = "done"
$ws.Range("E13").Value = "done"

# Flag the "view order" row (12) with a note about an outstanding issue,
# highlighted with a yellow fill.
$ws.Range("E12").Value = "issue ,can not populate order data into db"
$ws.Range("E12").Interior.Color = 65535

# Leave the cursor where the author last left it.
$null = $ws.Range("D10").Select()
